# ---------------------------------------------------------------------------
# Applies two changes described by the commit's OOXML diff:
#
# 1. The table on slide 16 (the 3rd shape, a graphicFrame/table) switches
#    from the deck's custom "Table_0" style to the built-in table style
#    {166697A4-CA72-4192-8963-CF13912B0AE0}.
#
# 2. The presentation's theme colour palette changes from the "Integral"
#    palette to the standard "Office" palette (dk2/lt2/accent1-6/hlink/
#    folHlink all change; dk1/lt1 stay black/white in both palettes).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
if (-not $tableShape.HasTable) {
    # Defensive fallback: locate the table shape by scanning the slide.
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        if ($slide.Shapes.Item($j).HasTable) {
            $tableShape = $slide.Shapes.Item($j)
            break
        }
    }
}
$tableShape.Table.ApplyStyle("{166697A4-CA72-4192-8963-CF13912B0AE0}")

# --- 2. Theme colours: Integral -> Office ---------------------------------
# Colors() index order is dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
